$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the "Attr_1" values from "No" to "Str" for Giant, Human and Rootwalker rows
$ws.Range("D5").Value = "Str"
$ws.Range("D7").Value = "Str"
$ws.Range("D11").Value = "Str"

# Best-fit the first few columns (A, B, D) so their contents are fully visible
$ws.Columns.Item(1).ColumnWidth = 14.16666667
$ws.Columns.Item(2).ColumnWidth = 10.5
$ws.Columns.Item(4).ColumnWidth = 11.5

# Move the active selection to H30, matching the saved cursor position
$ws.Range("H30").Select() | Out-Null
